# "Categories and features are updated"
#
# The "research" sheet's holdout-ratio table is expanded from a single
# category (50% holdout) into three categories (50%, 33%, 25% holdout),
# with the HOLDOUT RATIO column switched from text ("0.5") to a real
# number, and new DESCRIPTION strings added for the two new rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("research")

# Row 3 and 4 are brand new rows - duplicate row 2's formatting (borders,
# number format, etc.) onto them first, then overwrite the values below.
$ws.Range("A2:D2").Copy() | Out-Null
$ws.Range("A3:D4").PasteSpecial(-4122) | Out-Null

# Row 2: category 1 - 50% holdout ratio (now a number, not text)
$ws.Range("B2").Value = 0.5
$ws.Range("D2").Value = "50% of the tracks are hidden, seeds are randomized"

# Row 3 (new): category 2 - 33% holdout ratio
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = 0.33
$ws.Range("C3").Value = $true
$ws.Range("D3").Value = "33% of the tracks are hidden, seeds are randomized"

# Row 4 (new): category 3 - 25% holdout ratio
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = 0.25
$ws.Range("C4").Value = $true
$ws.Range("D4").Value = "25% of the tracks are hidden, seeds are randomized"

# Author's cursor ended up on D4 after the edit
$ws.Range("D4").Select() | Out-Null
